# feat: add id on excel
#
# The sheet had a manually numbered "id" column (A5:A35, values 1..33)
# which is being cleared out (presumably to be replaced by an auto id /
# formula at a later point) while keeping the existing cell formatting.
# The view/selection state is also updated to reflect where the user was
# working (scrolled down one row, with A5:A19 selected and A19 active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the numeric id values in column A (rows 5 through 35) but keep the
# cell styling intact.
$ws.Range("A5:A35").ClearContents()

# Update the view: scroll so row 2 is the top-left visible row, and select
# A5:A19 with A19 as the active cell.
$ws.Range("A5:A19").Select()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("A19").Activate()
